$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 133, shifting existing rows 133:141 down to 134:142.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new data record.
$ws.Cells.Item(133, 1).Value = 5
$ws.Cells.Item(133, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(133, 3).Value = "Maule"
$ws.Cells.Item(133, 4).Value = 44461
$ws.Cells.Item(133, 5).Value = 7
$ws.Cells.Item(133, 6).Value = 100112008
$ws.Cells.Item(133, 7).Value = "Coliflor"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 3000
$ws.Cells.Item(133, 11).Value = 600
$ws.Cells.Item(133, 12).Value = 600
$ws.Cells.Item(133, 13).Value = 600
$ws.Cells.Item(133, 14).Value = "`$/unidad"
$ws.Cells.Item(133, 15).Value = "Región del Maule"
$ws.Cells.Item(133, 16).Value = 600
$ws.Cells.Item(133, 17).Value = 1
$ws.Cells.Item(133, 18).Value = "Hortaliza"

# Match the date-column number format used by the rest of column D.
$ws.Cells.Item(133, 4).NumberFormat = $ws.Cells.Item(134, 4).NumberFormat
